# Refined metadata to be additional tab
#
# 1. Update the F-column (time_taken) timestamps on the existing "data" sheet.
# 2. Add a new "metadata" sheet (after "data") describing the panelapp query
#    that produced this workbook.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- 1. Refresh the time_taken timestamps on the "data" sheet -------------
$timestamps = @(
  "2021-10-05 14:34:25.888722",
  "2021-10-05 14:34:25.888730",
  "2021-10-05 14:34:25.888733",
  "2021-10-05 14:34:25.888736",
  "2021-10-05 14:34:25.888739",
  "2021-10-05 14:34:25.888742",
  "2021-10-05 14:34:25.888744",
  "2021-10-05 14:34:25.888747",
  "2021-10-05 14:34:25.888750",
  "2021-10-05 14:34:25.888752",
  "2021-10-05 14:34:25.888755",
  "2021-10-05 14:34:25.888757",
  "2021-10-05 14:34:25.888760",
  "2021-10-05 14:34:25.888762",
  "2021-10-05 14:34:25.888765",
  "2021-10-05 14:34:25.888767",
  "2021-10-05 14:34:25.888770",
  "2021-10-05 14:34:25.888773",
  "2021-10-05 14:34:25.888775",
  "2021-10-05 14:34:25.888778"
)
for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = $i + 2
  $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- 2. Add the "metadata" sheet, placed right after "data" ---------------
$ws = $wb.Worksheets.Add($null, $data)
$ws.Name = "metadata"

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Lymphoedema_nonsyndromic"
$ws.Range("C2").Value = 133
$ws.Range("D2").Value = "'0.26"
$ws.Range("E2").Value = "2021-07-05T07:37:40.196084Z"
$ws.Range("F2").Value = "2021-10-05 14:34:25.885033"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/133/?format=json"

# Match the bold/centered/bordered header look already used on "data"
# (copy the formatting of the matching header cells rather than re-building
# fonts/borders from scratch, so the existing style is reused).
$data.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Keep "data" as the active/selected sheet, as it was before this edit.
$data.Activate()
